$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking strings are not converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.037.89"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "1.886.78"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "329.98"
$ws.Range("E5").Value = "  -3.42%  "

$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -4.07%  "

$ws.Range("E8").Value = "  -0.87%  "

$ws.Range("D9").Value = "47.66"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").Value = "0.07950"
$ws.Range("E10").Value = "  -3.64%  "

$ws.Range("D11").Value = "0.9970"
$ws.Range("E11").Value = "  -4.21%  "

$ws.Range("D12").Value = "21.77"
$ws.Range("E12").Value = "  -4.05%  "

$ws.Range("D13").Value = "1.891.16"
$ws.Range("E13").Value = "  -3.26%  "

$ws.Range("D14").Value = "5.913"
$ws.Range("E14").Value = "  -4.39%  "

$ws.Range("D15").Value = "7.067"
$ws.Range("E15").Value = "  -4.74%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").Value = "88.52"
$ws.Range("E17").Value = "  -4.27%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.06553"
$ws.Range("E18").Value = "  -1.79%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.00001024"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("E20").Value = "  -3.56%  "

$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").Value = "29.034.59"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").Value = "5.434"
$ws.Range("E23").Value = "  -3.09%  "

$ws.Range("D24").Value = "11.42"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("E25").Value = "  -3.65%  "

$ws.Range("D26").Value = "2.133.41"
$ws.Range("E26").Value = "  -2.38%  "

$ws.Range("D27").Value = "155.94"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("E28").Value = "  -3.05%  "

$ws.Range("D29").Value = "2.080"
$ws.Range("E29").Value = "  -5.12%  "

$ws.Range("D30").Value = "5.498"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").Value = "117.35"
$ws.Range("E31").Value = "  -4.09%  "

$ws.Range("D32").Value = "1.039"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").Value = "0.09309"
$ws.Range("E33").Value = "  -3.73%  "

$ws.Range("E34").Value = "  -4.04%  "

$ws.Range("D35").Value = "3.521"
$ws.Range("E35").Value = "  -4.37%  "

$ws.Range("D36").Value = "5.286"
$ws.Range("E36").Value = "  -3.79%  "

$ws.Range("E37").Value = "  -3.78%  "

$ws.Range("D38").Value = "0.02224"
$ws.Range("E38").Value = "  -4.05%  "

$ws.Range("D39").Value = "8.332"
$ws.Range("E39").Value = "  -3.42%  "

$ws.Range("D40").Value = "1.171"
$ws.Range("E40").Value = "  -1.84%  "

$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "0.5778"
$ws.Range("E42").Value = "  -5.42%  "

$ws.Range("D43").Value = "0.1822"
$ws.Range("E43").Value = "  -4.36%  "

$ws.Range("D44").Value = "10.08"
$ws.Range("E44").Value = "  -5.85%  "

$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("D46").Value = "0.07526"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").Value = "2.276"
$ws.Range("E47").Value = "  -2.83%  "

$ws.Range("E48").Value = "  -4.38%  "

$ws.Range("D49").Value = "0.5453"
$ws.Range("E49").Value = "  -4.68%  "

$ws.Range("D50").Value = "1.900"
$ws.Range("E50").Value = "  -5.06%  "

$ws.Range("D51").Value = "111.43"
$ws.Range("E51").Value = "  -2.74%  "
